$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Remis'
$ws.Range("A3").Value = 'Remis'
$ws.Range("A6").Value = 'Górnik Zabrze'
$ws.Range("A9").Value = 'Raków Częstochowa'
$ws.Range("A17").Value = 'Stal Mielec'
$ws.Range("A18").Value = 'Wisła Płock'
$ws.Range("A20").Value = 'Remis'
$ws.Range("A21").Value = 'Piast Gliwice'
$ws.Range("A24").Value = 'Remis'
$ws.Range("A26").Value = 'Radomiak Radom'
$ws.Range("A29").Value = 'Remis'
$ws.Range("A34").Value = 'Cracovia'
$ws.Range("A35").Value = 'Pogoń Szczecin'
$ws.Range("A37").Value = 'Śląsk Wrocław'
$ws.Range("A39").Value = 'Korona Kielce'
$ws.Range("A40").Value = 'Śląsk Wrocław'
$ws.Range("A45").Value = 'Raków Częstochowa'
$ws.Range("A48").Value = 'Miedź Legnica'
$ws.Range("A49").Value = 'Remis'
$ws.Range("A50").Value = 'Lech Poznań'
$ws.Range("A53").Value = 'Warta Poznań'
$ws.Range("A58").Value = 'Remis'
$ws.Range("A59").Value = 'Jagielonia Białystok'
$ws.Range("A60").Value = 'Remis'
$ws.Range("A62").Value = 'Stal Mielec'
$ws.Range("A63").Value = 'Widzew Łódź'
$ws.Range("A66").Value = 'Piast Gliwice'
$ws.Range("A67").Value = 'Remis'
$ws.Range("A74").Value = 'Stal Mielec'
$ws.Range("A75").Value = 'Korona Kielce'
$ws.Range("A76").Value = 'Remis'
$ws.Range("A78").Value = 'Wisła Płock'
$ws.Range("A79").Value = 'Raków Częstochowa'
$ws.Range("A80").Value = 'Zagłębie Lubin'
$ws.Range("A81").Value = 'Widzew Łódź'
$ws.Range("A82").Value = 'Remis'
$ws.Range("A84").Value = 'Śląsk Wrocław'
$ws.Range("A86").Value = 'Górnik Zabrze'
$ws.Range("A90").Value = 'Widzew Łódź'
$ws.Range("A91").Value = 'Remis'
$ws.Range("A96").Value = 'Pogoń Szczecin'
$ws.Range("A97").Value = 'Remis'
$ws.Range("A99").Value = 'Remis'
$ws.Range("A100").Value = 'Śląsk Wrocław'
$ws.Range("A101").Value = 'Lechia Gdańsk'
$ws.Range("A102").Value = 'Piast Gliwice'
$ws.Range("A103").Value = 'Wisła Płock'
$ws.Range("A106").Value = 'Remis'
$ws.Range("A107").Value = 'Raków Częstochowa'
$ws.Range("A108").Value = 'Stal Mielec'
$ws.Range("A109").Value = 'Śląsk Wrocław'
$ws.Range("A112").Value = 'Miedź Legnica'
$ws.Range("A113").Value = 'Górnik Zabrze'
$ws.Range("A116").Value = 'Remis'
$ws.Range("A118").Value = 'Legia Warszawa'
$ws.Range("A121").Value = 'Lechia Gdańsk'
$ws.Range("A122").Value = 'Legia Warszawa'
$ws.Range("A123").Value = 'Korona Kielce'
$ws.Range("A124").Value = 'Remis'
$ws.Range("A125").Value = 'Warta Poznań'
$ws.Range("A126").Value = 'Widzew Łódź'
$ws.Range("A127").Value = 'Śląsk Wrocław'
$ws.Range("A128").Value = 'Legia Warszawa'
$ws.Range("A130").Value = 'Piast Gliwice'
$ws.Range("A132").Value = 'Lechia Gdańsk'
$ws.Range("A135").Value = 'Warta Poznań'
$ws.Range("A137").Value = 'Remis'
$ws.Range("A138").Value = 'Remis'
$ws.Range("A146").Value = 'Lech Poznań'
$ws.Range("A148").Value = 'Korona Kielce'
$ws.Range("A153").Value = 'Wisła Płock'
$ws.Range("A154").Value = 'Remis'
$ws.Range("A155").Value = 'Górnik Zabrze'
$ws.Range("A156").Value = 'Piast Gliwice'
$ws.Range("A158").Value = 'Legia Warszawa'
$ws.Range("A159").Value = 'Miedź Legnica'
$ws.Range("A163").Value = 'Zagłębie Lubin'
$ws.Range("A167").Value = 'Górnik Zabrze'
$ws.Range("A170").Value = 'Piast Gliwice'
$ws.Range("A171").Value = 'Jagielonia Białystok'
$ws.Range("A172").Value = 'Wisła Płock'
$ws.Range("A174").Value = 'Zagłębie Lubin'
$ws.Range("A180").Value = 'Lech Poznań'
$ws.Range("A181").Value = 'Korona Kielce'
$ws.Range("A183").Value = 'Remis'
$ws.Range("A184").Value = 'Remis'
$ws.Range("A185").Value = 'Lech Poznań'
$ws.Range("A186").Value = 'Remis'
$ws.Range("A190").Value = 'Remis'
$ws.Range("A191").Value = 'Piast Gliwice'
$ws.Range("A197").Value = 'Remis'
$ws.Range("A198").Value = 'Pogoń Szczecin'
$ws.Range("A201").Value = 'Wisła Płock'
$ws.Range("A202").Value = 'Remis'
$ws.Range("A203").Value = 'Jagielonia Białystok'
$ws.Range("A204").Value = 'Remis'
$ws.Range("A206").Value = 'Zagłębie Lubin'
$ws.Range("A207").Value = 'Piast Gliwice'
$ws.Range("A208").Value = 'Remis'
$ws.Range("A209").Value = 'Piast Gliwice'
$ws.Range("A210").Value = 'Jagielonia Białystok'
$ws.Range("A211").Value = 'Remis'
$ws.Range("A220").Value = 'Górnik Zabrze'
$ws.Range("A221").Value = 'Pogoń Szczecin'
$ws.Range("A222").Value = 'Remis'
$ws.Range("A224").Value = 'Warta Poznań'
$ws.Range("A227").Value = 'Cracovia'
$ws.Range("A229").Value = 'Remis'
$ws.Range("A231").Value = 'Pogoń Szczecin'
$ws.Range("A235").Value = 'Remis'
$ws.Range("A236").Value = 'Remis'
$ws.Range("A243").Value = 'Wisła Płock'
$ws.Range("A244").Value = 'Remis'
$ws.Range("A245").Value = 'Cracovia'
$ws.Range("A246").Value = 'Wisła Płock'
$ws.Range("A248").Value = 'Jagielonia Białystok'
$ws.Range("A252").Value = 'Stal Mielec'
$ws.Range("A253").Value = 'Remis'
$ws.Range("A254").Value = 'Zagłębie Lubin'
$ws.Range("A255").Value = 'Remis'
$ws.Range("A258").Value = 'Pogoń Szczecin'
$ws.Range("A260").Value = 'Legia Warszawa'
$ws.Range("A261").Value = 'Piast Gliwice'
$ws.Range("A262").Value = 'Jagielonia Białystok'
$ws.Range("A265").Value = 'Remis'
$ws.Range("A269").Value = 'Raków Częstochowa'
$ws.Range("A270").Value = 'Stal Mielec'
$ws.Range("A271").Value = 'Remis'
$ws.Range("A274").Value = 'Lech Poznań'
$ws.Range("A278").Value = 'Remis'
$ws.Range("A279").Value = 'Piast Gliwice'
$ws.Range("A280").Value = 'Stal Mielec'
$ws.Range("A285").Value = 'Raków Częstochowa'
$ws.Range("A289").Value = 'Wisła Płock'
$ws.Range("A291").Value = 'Zagłębie Lubin'
$ws.Range("A293").Value = 'Remis'
$ws.Range("A294").Value = 'Górnik Zabrze'
$ws.Range("A295").Value = 'Remis'
$ws.Range("A298").Value = 'Śląsk Wrocław'
$ws.Range("A299").Value = 'Wisła Płock'
$ws.Range("A305").Value = 'Raków Częstochowa'
$ws.Range("A307").Value = 'Widzew Łódź'
